$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "05:35:00"
$ws.Range("L3").Value = 1.28
$ws.Range("R4").Value = 1.18
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 3.8
$ws.Range("O5").Value = 1.28
$ws.Range("P5").Value = 1.98
$ws.Range("R5").Value = 1.37
$ws.Range("S5").Value = 3.1
$ws.Range("T5").Value = 1.86
$ws.Range("U5").Value = 1.92
$ws.Range("V5").Value = 1.16
$ws.Range("W5").Value = 2.62
$ws.Range("X5").Value = 17
$ws.Range("Y5").Value = 22
$ws.Range("Z5").Value = 55
$ws.Range("AA5").Value = 210
$ws.Range("AB5").Value = 8.800000000000001
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 25
$ws.Range("AE5").Value = 110
$ws.Range("AF5").Value = 10.5
$ws.Range("AG5").Value = 10.5
$ws.Range("AH5").Value = 23
$ws.Range("AI5").Value = 100
$ws.Range("AJ5").Value = 17
$ws.Range("AK5").Value = 18.5
$ws.Range("AL5").Value = 38
$ws.Range("AM5").Value = 140
$ws.Range("AN5").Value = 10
$ws.Range("AO5").Value = 130
$ws.Range("I6").Value = 2.34
$ws.Range("J6").Value = 2.58
$ws.Range("Q7").Value = 1.83
$ws.Range("F8").Value = 4.6
$ws.Range("I8").Value = 1.3
$ws.Range("J8").Value = 6.6
$ws.Range("G9").Value = 2.38
$ws.Range("P9").Value = 1.28
$ws.Range("G10").Value = 3.4
$ws.Range("H10").Value = 2.86
$ws.Range("I10").Value = 3.1
$ws.Range("J10").Value = 2.82
$ws.Range("K10").Value = 3.05
$ws.Range("P10").Value = 1.49
$ws.Range("Q10").Value = 2.8
$ws.Range("G12").Value = 1.45
$ws.Range("Q12").Value = 1.7
$ws.Range("R12").Value = 1.52
$ws.Range("U12").Value = 1.97
$ws.Range("AD13").Value = 19.5
$ws.Range("AN13").Value = 32
$ws.Range("AH14").Value = 26
